# "Upload Leave Card 12/27/2023 4:01 PM"
#
# Rolls the leave-card workbook forward one period:
#  - "2018 LEAVE CREDITS" sheet: shifts the monthly PERIOD dates (rows 77-95)
#    from month-start to month-end, and backfills the EARNED column for the
#    period that had just closed (rows 80-86), which ripples into the
#    BALANCE formulas in row 9 (E9 / I9) automatically.
#  - "2017 LEAVE BALANCE" sheet: records a new VL leave transaction in row 51
#    (PERIOD 11/1/2023, PARTICULARS "VL(6-0-0)", 6 days charged, REMARKS the
#    leave dates), which ripples into the BALANCE formula in row 9 (E9).
#  - Leaves "2017 LEAVE BALANCE" as the active sheet/tab, matching the last
#    thing the author touched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "2018 LEAVE CREDITS" sheet
# ---------------------------------------------------------------------
$ws2018 = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws2018.Activate() | Out-Null

# Monthly PERIOD date column (A) shifts from the 1st of each month to the
# last day of each month, rows 77-95.
$periodDates = @{
    77 = 44957; 78 = 44985; 79 = 45016; 80 = 45046; 81 = 45077;
    82 = 45107; 83 = 45138; 84 = 45169; 85 = 45199; 86 = 45230;
    87 = 45260; 88 = 45291; 89 = 45322; 90 = 45351; 91 = 45382;
    92 = 45412; 93 = 45443; 94 = 45473; 95 = 45504
}
foreach ($row in $periodDates.Keys) {
    $ws2018.Range("A$row").Value = $periodDates[$row]
}

# Newly-earned leave credits (1.25 days) posted for the periods that just
# closed out; column G recalculates off of column C automatically.
foreach ($row in 80..86) {
    $ws2018.Range("C$row").Value = 1.25
}

# ---------------------------------------------------------------------
# "2017 LEAVE BALANCE" sheet
# ---------------------------------------------------------------------
$ws2017 = $wb.Worksheets.Item("2017 LEAVE BALANCE")
$ws2017.Activate() | Out-Null

$ws2017.Range("A51").Value = 45231
$ws2017.Range("B51").Value = "VL(6-0-0)"
$ws2017.Range("D51").Value = 6
$ws2017.Range("K51").Value = "11/22-24,27-29/2023"

# ---------------------------------------------------------------------
# View state: leave the active selection on each sheet roughly where the
# author left it, with "2017 LEAVE BALANCE" as the active tab. The split
# pane's bottom-left half (containing the data rows) is where the author's
# selection ultimately rests on each sheet.
# ---------------------------------------------------------------------
$ws2018.Range("I9").Select() | Out-Null
$ws2018.Range("G89").Select() | Out-Null

$ws2017.Activate() | Out-Null
$ws2017.Range("I9").Select() | Out-Null
$ws2017.Range("K51").Select() | Out-Null
